$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Fecha (D) and Volumen (M) values per row, reflecting the
# weekly re-shuffle of dates/volumes across rows 2,3,4,6,7,9,10.
$updates = @(
    @{ Row = 2;  D = 44313; M = 120 },
    @{ Row = 3;  D = 44306; M = 80 },
    @{ Row = 4;  D = 44316; M = 120 },
    @{ Row = 6;  D = 44327; M = 60 },
    @{ Row = 7;  D = 44322; M = 60 },
    @{ Row = 9;  D = 44302; M = 80 },
    @{ Row = 10; D = 44323; M = 80 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D   # Column D = Fecha
    $ws.Cells.Item($u.Row, 13).Value = $u.M  # Column M = Volumen
}
